# Update 1000 runs ascended data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Optimizer | 1000 run - Ascended")

# New raw C/D column inputs (rows 3-13). Column E/L/M/Q/R recompute
# automatically from the existing formulas already on the sheet.
$newValues = @{
    3  = @(281733, 83877)
    4  = @(281733, 355794)
    5  = @(281733, 693264)
    6  = @(281733, 1067934)
    7  = @(281733, 1479804)
    8  = @(281733, 1928874)
    9  = @(281733, 2415144)
    10 = @(281733, 2938614)
    11 = @(281733, 3499284)
    12 = @(281733, 4097154)
    13 = @(281733, 4712149)
}

foreach ($row in $newValues.Keys) {
    $pair = $newValues[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}

# The author left this sheet active/selected (cell H21) when the workbook
# was last saved, moving the selection away from "200 run - Descended".
$ws.Activate()
$ws.Range("H21").Select()
